# FALECPV-RecAnulados.xlsx: remove the "AUT S.R.I." column from the
# "RECIBOS ANULADOS" report table (Hoja1). Deleting the whole column lets
# Excel shift CANTIDAD/SUBTOTAL/DESCUENTO/TOTAL/PAGO/ENTREGA/CAMBIO left by
# one, shrink the merged title range, and drop the now-unused shared string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("F").EntireColumn.Delete()

# Matches the post-edit selection recorded in the saved sheet view.
$ws.Range("J11").Select()
